$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values (changed in the update)
$ws.Range("B74").Value = 21699
$ws.Range("C74").Value = -1225
$ws.Range("D74").Value = 20473
$ws.Range("E74").Value = 10577
$ws.Range("F74").Value = -2891
$ws.Range("G74").Value = 178
$ws.Range("H74").Value = 7183
$ws.Range("I74").Value = 10886
$ws.Range("J74").Value = 20
$ws.Range("K74").Value = -3683

# Add new row 75 with the new quarter
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = 22996
$ws.Range("C75").Value = 178
$ws.Range("D75").Value = 23174
$ws.Range("E75").Value = 10899
$ws.Range("F75").Value = -2841
$ws.Range("G75").Value = -207
$ws.Range("H75").Value = 9227
$ws.Range("I75").Value = 9562
$ws.Range("J75").Value = 323
$ws.Range("K75").Value = -11
